$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LR-pair data table (rows 2-11, columns A-T), replacing/extending the
# previous rows 2-6. Values reflect the refreshed TPM computation: every
# sending cluster now appears once per target cluster (FAPs and MuSCs).
$data = @(
    @("ECs","L1cam","Cntn1","FAPs",3,1,6.713252999999999,20.139759,0.3101840064655811,0.3231642354899327,2,0.6666666666666666,0.07031566666666667,0.210947,0.3500709860127268,0.446885632088942,0.4720468601969999,4.248421741773,0.108586420988784,0.1444174536454583),
    @("ECs","L1cam","Cntn1","MuSCs",3,1,6.713252999999999,20.139759,0.3101840064655811,0.3231642354899327,2,1,0.1305455,0.261091,0.6499290139872732,0.5531143679110581,0.8763849695114999,5.258309817069,0.2015975854767971,0.1787467818444744),
    @("FAPs","L1cam","Cntn1","FAPs",3,1,0.3343473333333333,1.003042,0.0154484264788496,0.01609489473505086,2,0.6666666666666666,0.07031566666666667,0.210947,0.3500709860127268,0.446885632088942,0.02350985564155555,0.211588700774,0.005408045889795996,0.007192577207078187),
    @("FAPs","L1cam","Cntn1","MuSCs",3,1,0.3343473333333333,1.003042,0.0154484264788496,0.01609489473505086,2,1,0.1305455,0.261091,0.6499290139872732,0.5531143679110581,0.04364753980366667,0.261885238822,0.0100403805890536,0.008902317527972671),
    @("Inflammatory-Mac","L1cam","Cntn1","FAPs",3,1,6.661784666666667,19.985354,0.3078059262949933,0.3206866401135023,2,0.6666666666666666,0.07031566666666667,0.210947,0.3500709860127268,0.446885632088942,0.4684278300264444,4.215850470238,0.107753924118649,0.1433102518696015),
    @("Inflammatory-Mac","L1cam","Cntn1","MuSCs",3,1,6.661784666666667,19.985354,0.3078059262949933,0.3206866401135023,2,1,0.1305455,0.261091,0.6499290139872732,0.5531143679110581,0.8696660102023334,5.217996061214,0.2000520021763443,0.1773763882439008),
    @("MuSCs","L1cam","Cntn1","FAPs",2,1,2.607918,5.215835999999999,0.1204981331366039,0.08369373503331734,2,0.6666666666666666,0.07031566666666667,0.210947,0.3500709860127268,0.446885632088942,0.183377492782,1.100264956692,0.04218290027982376,0.03740152768224845),
    @("MuSCs","L1cam","Cntn1","MuSCs",2,1,2.607918,5.215835999999999,0.1204981331366039,0.08369373503331734,2,1,0.1305455,0.261091,0.6499290139872732,0.5531143679110581,0.340451959269,1.361807837076,0.07831523285678013,0.0462922073510689),
    @("Resolving-Mac","L1cam","Cntn1","FAPs",3,1,5.325505333333333,15.976516,0.2460635076239721,0.2563604946281968,2,0.6666666666666666,0.07031566666666667,0.210947,0.3500709860127268,0.446885632088942,0.3744664578502222,3.370198120652,0.08613969473567405,0.1145638216845555),
    @("Resolving-Mac","L1cam","Cntn1","MuSCs",3,1,5.325505333333333,15.976516,0.2460635076239721,0.2563604946281968,2,1,0.1305455,0.261091,0.6499290139872732,0.5531143679110581,0.6952207564926667,4.171324538956,0.1599238128882981,0.1417966729436413)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item(2 + $r, 1 + $c).Value = $rowValues[$c]
    }
}

Write-Output ("Updated range A2:" + $ws.Cells.Item(1 + $data.Count, $data[0].Count).Address($false, $false))
